$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(46).Insert()

$ws.Range("A46").Value = 4
$ws.Range("B46").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C46").Value = "Los Lagos"
$ws.Range("D46").Value = 44930
$ws.Range("E46").Value = 10
$ws.Range("F46").Value = 100112026
$ws.Range("G46").Value = "Haba"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 15
$ws.Range("K46").Value = 30000
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = 30000
$ws.Range("N46").Value = "`$/saco 25 kilos"
$ws.Range("O46").Value = "Región de La Araucanía"
$ws.Range("P46").Value = 1200
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"
